$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Squad Total" row (row 35): remove the label in A35 and all
# summary formulas/values in B35:U35, while keeping cell formatting intact.
$ws.Range("A35:U35").ClearContents()

# Update the view: select the full row 35 (matches the saved selection in
# the edited workbook).
$ws.Rows("35:35").Select()
